$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update sprint dates (row 3 / row 4)
$ws.Range("D3").Value = 43656
$ws.Range("C4").Value = 43655

# Update row 5 (MAX Lane) contents
$ws.Range("D5").Value = "PHP 7 Regression"
$ws.Range("E5").Value = "PHP 7 Launch"

# Update row 6 contents
$ws.Range("E6").Value = "Business Case / SIR Updates"
$ws.Range("F6").Value = "SIL Updates"

# Update row 7 contents
$ws.Range("D7").Value = ""
$ws.Range("F7").Value = "UAT Findings"

# Update row 8 contents
$ws.Range("E8").Value = "UAT Findings"

# Update selection
$ws.Range("C5").Select()
